$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 217-218; this shifts all existing data (rows 217-350)
# down to rows 219-352, and grows the used range from A1:R350 to A1:R352.
$ws.Rows("217:218").Insert()

# Populate the newly inserted row 217 with this week's "Primera" quality data.
$ws.Range("A217").Value = 3
$ws.Range("B217").Value = "Femacal de La Calera"
$ws.Range("C217").Value = "Coquimbo"
$ws.Range("D217").Value = 44438
$ws.Range("E217").Value = 5
$ws.Range("F217").Value = 100114014
$ws.Range("G217").Value = "Betarraga"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 3100
$ws.Range("K217").Value = 600
$ws.Range("L217").Value = 650
$ws.Range("M217").Value = 626
$ws.Range("N217").Value = "$/paquete 4 unidades"
$ws.Range("O217").Value = "Provincia de Quillota"
$ws.Range("P217").Value = 156
$ws.Range("Q217").Value = 4
$ws.Range("R217").Value = "Hortaliza"

# Populate the newly inserted row 218 with this week's "Segunda" quality data.
$ws.Range("A218").Value = 3
$ws.Range("B218").Value = "Femacal de La Calera"
$ws.Range("C218").Value = "Coquimbo"
$ws.Range("D218").Value = 44438
$ws.Range("E218").Value = 5
$ws.Range("F218").Value = 100114014
$ws.Range("G218").Value = "Betarraga"
$ws.Range("H218").Value = "Sin especificar"
$ws.Range("I218").Value = "Segunda"
$ws.Range("J218").Value = 1800
$ws.Range("K218").Value = 400
$ws.Range("L218").Value = 400
$ws.Range("M218").Value = 400
$ws.Range("N218").Value = "$/paquete 4 unidades"
$ws.Range("O218").Value = "Provincia de Quillota"
$ws.Range("P218").Value = 100
$ws.Range("Q218").Value = 4
$ws.Range("R218").Value = "Hortaliza"
